# Applies the "lookup.xlsx" edit described in the commit:
#   - Corrected label for nonmcbaddebt field (thank you Ken Michelson for
#     finding this bug)
#   - Added labels for disabled uncompensated care variables should they
#     eventually be enabled.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Lookup Table")
$ws2 = $wb.Worksheets.Item("Type and Label")

# ---------------------------------------------------------------------
# 1) Fix the mislabeled "nonmcbaddebt" row (rec/type/label row 22 on the
#    "Type and Label" sheet).
# ---------------------------------------------------------------------
$ws2.Range("C22").Value = "non-medicare bad debt expense (2010 format only)"

# ---------------------------------------------------------------------
# 2) Insert five new label rows (35-39) for the previously-unlabeled,
#    disabled uncompensated-care / bad-debt variables, pushing the old
#    row 35 (chguccare) down to row 40.
# ---------------------------------------------------------------------
$ws2.Range("A35:A39").EntireRow.Insert()

$newRows = @(
    ,@("costinitchcare", "dollar_flow", "cost of patients approved for charity care and uninsured discounts (2010 format only)")
    ,@("costchcare",     "dollar_flow", "cost of charity care (2010 format only)")
    ,@("totbaddebt",     "dollar_flow", "total bad debt expense (2010 format only)")
    ,@("mcbaddebt",      "dollar_flow", "medicare reimbursable bad debts (2010 format only)")
    ,@("baddebt",        "dollar_flow", "cost of non-Medicare and non-reimbursable Medicare bad debt expense (2010 format only)")
)

$r = 35
foreach ($row in $newRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) View-state bookkeeping: "Type and Label" becomes the active/front
#    tab, with its own scroll/selection; "Lookup Table" keeps its own
#    scroll/selection but is no longer the active tab.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A35:A39").Select()

$ws2.Activate()
$ws2.Range("C36").Select()
